# IPA new script implementation
# Adds two new test-case rows (IPA5 / IPA6) to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Row heights for the two new rows (105 / 120 points respectively)
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 105
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------------
# 2. Cell values — written in the exact order that reproduces the target
#    shared-string table ordering: IPA5, IPA6, OPQA-4205..., OPQA-4197...,
#    meta-data paragraph, terms-of-use paragraph.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "IPA5"
$ws.Range("A16").Value = "IPA6"

$ws.Range("B15").Value = "OPQA-4205||OPQA-4207||OPQA-4208||OPQA-4210||OPQA-4211"
$ws.Range("B16").Value = "OPQA-4197||OPQA-4199||OPQA-4215||OPQA-4216||OPQA-4201"

$ws.Range("C16").Value = "Verify that profile fly-out will display profile meta-data||Verify that profile fly-out provides access to the profile modal.||Verify that the profile fly-out should display the following user profile details, if available: a)First name b)Last Name c)Title d)Institution e)Country f)Photo||Verify that by clicking on any of the following fields (when present), will provide access to the profile modal. 1.Name 2.Institution 3. Country 4 .Title 5.Photo || Verify that profile fly-out provides access to the account setting modal"
$ws.Range("C15").Value = "Verify that the profile fly-out should contain link to terms of use||Verify that profile fly-out should contain link to privacy statement||Verify that the profile fly-out should contain link to app-specific feedback page||Verify that the profile fly-out should contain link to app-specific help page||Verify that the alternative profile fly-out should contain link to sign out of the platform. User returns to sign-in page."

$ws.Range("D15").Value = "Y"
$ws.Range("D16").Value = "Y"

# ---------------------------------------------------------------------------
# 3. Cell formatting — copy from existing rows with equivalent border /
#    wrap-text / fill styling so the new rows look consistent with the rest
#    of the table.
# ---------------------------------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B6").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D7").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("A8").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B11").Copy()
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E8").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. View state — scroll the sheet down a bit and move the active selection,
#    matching the author's on-screen state when the rows were added.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H15").Select()
